# Update maps/BloodVolumeControlv1.xlsx with component and flatmap id columns.
#
# The sheet originally had columns A:E (Name, Unit, Initial Value,
# Description/Physiological meaning, Relevant map feature). This change
# inserts a new "Component" column at D (pushing the old D/E to E/F) and
# appends a new "Flatmap ID" column at G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the new "Component" column before the old column D ---------
# Excel's column-insert copies the format of the column to the left (C),
# which is exactly the formatting the new column should carry.
$ws.Columns("D:D").Insert()

$ws.Range("D1").Value = "Component"

$componentRows = 3,4,5,6,7,8,10,11,12,13,14,15,16,17,18,19,20,22,23,24,25,26,27,28,29,30,31,32
foreach ($r in $componentRows) {
    $ws.Cells.Item($r, 4).Value = "main"
}

# --- 2. Append the new "Flatmap ID" column at G ----------------------------
$ws.Range("G1").Value = "Flatmap ID"

# Rows belonging to the "Circulation" / "GI tract" sections (3-20) got an
# explicit black font colour applied to the new Flatmap ID cells; rows in
# the "Kidney" section (22-32) were left with the default formatting.
$flatmapIds = @(
    @{ Row = 3;  Value = "bvc/ID-0000012"; Styled = $true },
    @{ Row = 4;  Value = "bvc/ID-0000009"; Styled = $true },
    @{ Row = 6;  Value = "bvc/ID-0000038"; Styled = $true },
    @{ Row = 7;  Value = "bvc/ID-0000010"; Styled = $true },
    @{ Row = 8;  Value = "bvc/ID-0000011"; Styled = $true },
    @{ Row = 10; Value = "bvc/ID-0000033"; Styled = $true },
    @{ Row = 11; Value = "bvc/ID-0000029"; Styled = $true },
    @{ Row = 12; Value = "bvc/ID-0000036"; Styled = $true },
    @{ Row = 13; Value = "bvc/ID-0000050"; Styled = $true },
    @{ Row = 14; Value = "bvc/ID-0000028"; Styled = $true },
    @{ Row = 15; Value = "bvc/ID-0000028"; Styled = $true },
    @{ Row = 16; Value = "bvc/ID-0000032"; Styled = $true },
    @{ Row = 17; Value = "bvc/ID-0000064"; Styled = $true },
    @{ Row = 18; Value = "bvc/ID-0000053"; Styled = $true },
    @{ Row = 19; Value = "bvc/ID-0000031"; Styled = $true },
    @{ Row = 20; Value = "bvc/ID-0000065"; Styled = $true },
    @{ Row = 22; Value = "bvc/ID-0000043"; Styled = $false },
    @{ Row = 23; Value = "bvc/ID-0000021"; Styled = $false },
    @{ Row = 24; Value = "bvc/ID-0000040"; Styled = $false },
    @{ Row = 25; Value = "bvc/ID-0000019"; Styled = $false },
    @{ Row = 26; Value = "bvc/ID-0000019"; Styled = $false },
    @{ Row = 27; Value = "bvc/ID-0000066"; Styled = $false },
    @{ Row = 28; Value = "bvc/ID-0000059"; Styled = $false },
    @{ Row = 29; Value = "bvc/ID-0000067"; Styled = $false },
    @{ Row = 30; Value = "bvc/ID-0000017"; Styled = $false },
    @{ Row = 31; Value = "bvc/ID-0000041"; Styled = $false },
    @{ Row = 32; Value = "bvc/ID-0000023"; Styled = $false }
)

foreach ($entry in $flatmapIds) {
    $cell = $ws.Cells.Item($entry.Row, 7)
    $cell.Value = $entry.Value
    if ($entry.Styled) {
        $cell.Font.Color = 0
    }
}
